$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 53
$ws_ALC.Range("H53").Value = 300.1
$ws_ALC.Range("I53").Value = 126.90909
$ws_ALC.Range("K53").Value = 126.90909
$ws_ALC.Range("M53").Value = 510.09091

# ALC row 55
$ws_ALC.Range("H55").Value = 269.4
$ws_ALC.Range("J55").Value = 441.6
$ws_ALC.Range("L55").Value = 441.6
$ws_ALC.Range("N55").Value = -869.6

# ALC row 86
$ws_ALC.Range("H86").Value = 37038790
$ws_ALC.Range("I86").Value = 45456452
$ws_ALC.Range("J86").Value = 1104.4
$ws_ALC.Range("K86").Value = 45456452
$ws_ALC.Range("L86").Value = 1104.4
$ws_ALC.Range("M86").Value = -45455329
$ws_ALC.Range("N86").Value = -3350.4

# ALC row 89
$ws_ALC.Range("H89").Value = 37038790
$ws_ALC.Range("I89").Value = 45456452
$ws_ALC.Range("J89").Value = 1104.4
$ws_ALC.Range("K89").Value = 227282260
$ws_ALC.Range("L89").Value = 5522
$ws_ALC.Range("M89").Value = -227276644
$ws_ALC.Range("N89").Value = -16754

# ALC row 111
$ws_ALC.Range("H111").Value = 13712
$ws_ALC.Range("J111").Value = 10998.5
$ws_ALC.Range("L111").Value = 32995.5
$ws_ALC.Range("N111").Value = -39129.5

# ALC row 112
$ws_ALC.Range("H112").Value = 2879.9412
$ws_ALC.Range("J112").Value = 2901.94
$ws_ALC.Range("L112").Value = 8705.82
$ws_ALC.Range("N112").Value = -10921.82

# ALC row 116
$ws_ALC.Range("H116").Value = 1026294.94
$ws_ALC.Range("J116").Value = 5330.375
$ws_ALC.Range("L116").Value = 5330.375
$ws_ALC.Range("N116").Value = -12214.375

# ALC row 132
$ws_ALC.Range("H132").Value = 173823.92
$ws_ALC.Range("I132").Value = 199152.47
$ws_ALC.Range("K132").Value = 597457.41
$ws_ALC.Range("M132").Value = -594927.41

# ALC row 138
$ws_ALC.Range("H138").Value = 5642.078
$ws_ALC.Range("J138").Value = 5688.8096
$ws_ALC.Range("L138").Value = 17066.4288
$ws_ALC.Range("N138").Value = -27346.4288

# ARM row 2
$ws_ARM.Range("H2").Value = 47329.727
$ws_ARM.Range("I2").Value = 168212.67
$ws_ARM.Range("J2").Value = 1998.625
$ws_ARM.Range("K2").Value = 168212.67
$ws_ARM.Range("L2").Value = 1998.625
$ws_ARM.Range("M2").Value = -168099.67
$ws_ARM.Range("N2").Value = -2224.625

# ARM row 32
$ws_ARM.Range("H32").Value = 5317.63
$ws_ARM.Range("I32").Value = 2666.854
$ws_ARM.Range("K32").Value = 2666.854
$ws_ARM.Range("M32").Value = -2379.854

# ARM row 61
$ws_ARM.Range("H61").Value = 5344.15
$ws_ARM.Range("I61").Value = 5093.5625
$ws_ARM.Range("J61").Value = 6346.5
$ws_ARM.Range("K61").Value = 5093.5625
$ws_ARM.Range("L61").Value = 6346.5
$ws_ARM.Range("M61").Value = -4881.5625
$ws_ARM.Range("N61").Value = -6770.5

# ARM row 74
$ws_ARM.Range("H74").Value = 10418841
$ws_ARM.Range("I74").Value = 12502050
$ws_ARM.Range("K74").Value = 12502050
$ws_ARM.Range("M74").Value = -12501176

# ARM row 77
$ws_ARM.Range("H77").Value = 10418841
$ws_ARM.Range("I77").Value = 12502050
$ws_ARM.Range("K77").Value = 62510250
$ws_ARM.Range("M77").Value = -62505882

# ARM row 110
$ws_ARM.Range("H110").Value = 2629
$ws_ARM.Range("I110").Value = 2569.7693
$ws_ARM.Range("J110").Value = 3399
$ws_ARM.Range("K110").Value = 2569.7693
$ws_ARM.Range("L110").Value = 3399
$ws_ARM.Range("M110").Value = -524.7692999999999
$ws_ARM.Range("N110").Value = -7489

# ARM row 116
$ws_ARM.Range("H116").Value = 47329.727
$ws_ARM.Range("I116").Value = 168212.67
$ws_ARM.Range("J116").Value = 1998.625
$ws_ARM.Range("K116").Value = 168212.67
$ws_ARM.Range("L116").Value = 1998.625
$ws_ARM.Range("M116").Value = -165918.67
$ws_ARM.Range("N116").Value = -6586.625

# ARM row 122
$ws_ARM.Range("H122").Value = 7308.8887
$ws_ARM.Range("I122").Value = 4757.6
$ws_ARM.Range("J122").Value = 10498
$ws_ARM.Range("K122").Value = 14272.8
$ws_ARM.Range("L122").Value = 31494
$ws_ARM.Range("M122").Value = -11822.8
$ws_ARM.Range("N122").Value = -36394

# ARM row 132
$ws_ARM.Range("H132").Value = 13629.679
$ws_ARM.Range("I132").Value = 16111.143
$ws_ARM.Range("K132").Value = 48333.429
$ws_ARM.Range("M132").Value = -45803.429

# ARM row 136
$ws_ARM.Range("H136").Value = 5344.15
$ws_ARM.Range("I136").Value = 5093.5625
$ws_ARM.Range("J136").Value = 6346.5
$ws_ARM.Range("K136").Value = 15280.6875
$ws_ARM.Range("L136").Value = 19039.5
$ws_ARM.Range("M136").Value = -12730.6875
$ws_ARM.Range("N136").Value = -24139.5

# BSM row 3
$ws_BSM.Range("H3").Value = 47329.727
$ws_BSM.Range("I3").Value = 168212.67
$ws_BSM.Range("J3").Value = 1998.625
$ws_BSM.Range("K3").Value = 168212.67
$ws_BSM.Range("L3").Value = 1998.625
$ws_BSM.Range("M3").Value = -168098.67
$ws_BSM.Range("N3").Value = -2226.625

# BSM row 20
$ws_BSM.Range("H20").Value = 4663.702
$ws_BSM.Range("I20").Value = 4050.3794
$ws_BSM.Range("J20").Value = 5651.8335
$ws_BSM.Range("K20").Value = 4050.3794
$ws_BSM.Range("L20").Value = 5651.8335
$ws_BSM.Range("M20").Value = -3803.3794
$ws_BSM.Range("N20").Value = -6145.8335

# BSM row 99
$ws_BSM.Range("H99").Value = 3800.682
$ws_BSM.Range("I99").Value = 4100.8823
$ws_BSM.Range("K99").Value = 4100.8823
$ws_BSM.Range("M99").Value = -2602.8823

# BSM row 134
$ws_BSM.Range("H134").Value = 5507.357
$ws_BSM.Range("I134").Value = 2945.818
$ws_BSM.Range("J134").Value = 14899.667
$ws_BSM.Range("K134").Value = 8837.454000000002
$ws_BSM.Range("L134").Value = 44699.001
$ws_BSM.Range("M134").Value = -6302.454000000002
$ws_BSM.Range("N134").Value = -49769.001

# CRP row 22
$ws_CRP.Range("H22").Value = 403.375
$ws_CRP.Range("I22").Value = 297.77777
$ws_CRP.Range("J22").Value = 539.1429000000001
$ws_CRP.Range("K22").Value = 297.77777
$ws_CRP.Range("L22").Value = 539.1429000000001
$ws_CRP.Range("M22").Value = 52.22223000000002
$ws_CRP.Range("N22").Value = -1239.1429

# CRP row 132
$ws_CRP.Range("H132").Value = 2195.1667
$ws_CRP.Range("I132").Value = 1977.9286
$ws_CRP.Range("K132").Value = 5933.7858
$ws_CRP.Range("M132").Value = -3403.7858

# CUL row 2
$ws_CUL.Range("H2").Value = 1430.8
$ws_CUL.Range("I2").Value = 20.571428
$ws_CUL.Range("J2").Value = 2664.75
$ws_CUL.Range("K2").Value = 123.428568
$ws_CUL.Range("L2").Value = 15988.5
$ws_CUL.Range("M2").Value = -10.42856800000001
$ws_CUL.Range("N2").Value = -16214.5

# CUL row 132
$ws_CUL.Range("H132").Value = 2491.0833
$ws_CUL.Range("J132").Value = 4632
$ws_CUL.Range("L132").Value = 41688
$ws_CUL.Range("N132").Value = -46748

# CUL row 139
$ws_CUL.Range("H139").Value = 2794.2104
$ws_CUL.Range("I139").Value = 1839.75
$ws_CUL.Range("J139").Value = 4430.4287
$ws_CUL.Range("K139").Value = 5519.25
$ws_CUL.Range("L139").Value = 13291.2861
$ws_CUL.Range("M139").Value = -379.25
$ws_CUL.Range("N139").Value = -23571.2861

# GSM row 95
$ws_GSM.Range("H95").Value = 32831.332
$ws_GSM.Range("J95").Value = 32831.332
$ws_GSM.Range("L95").Value = 32831.332
$ws_GSM.Range("N95").Value = -38323.332

# GSM row 102
$ws_GSM.Range("H102").Value = 56683220
$ws_GSM.Range("I102").Value = 170033330
$ws_GSM.Range("J102").Value = 8161.5
$ws_GSM.Range("K102").Value = 170033330
$ws_GSM.Range("L102").Value = 8161.5
$ws_GSM.Range("M102").Value = -170031708
$ws_GSM.Range("N102").Value = -11405.5

# GSM row 122
$ws_GSM.Range("H122").Value = 4502.121
$ws_GSM.Range("I122").Value = 2492.4546
$ws_GSM.Range("J122").Value = 5506.9546
$ws_GSM.Range("K122").Value = 7477.3638
$ws_GSM.Range("L122").Value = 16520.8638
$ws_GSM.Range("M122").Value = -5027.3638
$ws_GSM.Range("N122").Value = -21420.8638

# GSM row 126
$ws_GSM.Range("H126").Value = 4039.0476
$ws_GSM.Range("I126").Value = 2598.3667
$ws_GSM.Range("J126").Value = 7640.75
$ws_GSM.Range("K126").Value = 7795.1001
$ws_GSM.Range("L126").Value = 22922.25
$ws_GSM.Range("M126").Value = -5325.1001
$ws_GSM.Range("N126").Value = -27862.25

# GSM row 129
$ws_GSM.Range("H129").Value = 63997.5
$ws_GSM.Range("J129").Value = 63997.5
$ws_GSM.Range("L129").Value = 63997.5
$ws_GSM.Range("N129").Value = -73997.5

# GSM row 132
$ws_GSM.Range("H132").Value = 49534.316
$ws_GSM.Range("I132").Value = 57905.73
$ws_GSM.Range("J132").Value = 5285.4287
$ws_GSM.Range("K132").Value = 173717.19
$ws_GSM.Range("L132").Value = 15856.2861
$ws_GSM.Range("M132").Value = -171187.19
$ws_GSM.Range("N132").Value = -20916.2861

# LTW row 22
$ws_LTW.Range("H22").Value = 969.2174
$ws_LTW.Range("J22").Value = 1099.4286
$ws_LTW.Range("L22").Value = 1099.4286
$ws_LTW.Range("N22").Value = -1689.4286

# LTW row 27
$ws_LTW.Range("H27").Value = 969.2174
$ws_LTW.Range("J27").Value = 1099.4286
$ws_LTW.Range("L27").Value = 1099.4286
$ws_LTW.Range("N27").Value = -1313.4286

# LTW row 46
$ws_LTW.Range("H46").Value = 4544.8887
$ws_LTW.Range("J46").Value = 3898.3125
$ws_LTW.Range("L46").Value = 3898.3125
$ws_LTW.Range("N46").Value = -4274.3125

# LTW row 93
$ws_LTW.Range("H93").Value = 2690.3125
$ws_LTW.Range("I93").Value = 1911.875
$ws_LTW.Range("J93").Value = 3468.75
$ws_LTW.Range("K93").Value = 1911.875
$ws_LTW.Range("L93").Value = 3468.75
$ws_LTW.Range("M93").Value = -663.875
$ws_LTW.Range("N93").Value = -5964.75

# LTW row 132
$ws_LTW.Range("H132").Value = 5691.439
$ws_LTW.Range("I132").Value = 4942.04
$ws_LTW.Range("J132").Value = 8135.1304
$ws_LTW.Range("K132").Value = 14826.12
$ws_LTW.Range("L132").Value = 24405.3912
$ws_LTW.Range("M132").Value = -12296.12
$ws_LTW.Range("N132").Value = -29465.3912

# WVR row 88
$ws_WVR.Range("H88").Value = 24495.5
$ws_WVR.Range("J88").Value = 24495.5
$ws_WVR.Range("L88").Value = 24495.5
$ws_WVR.Range("N88").Value = -25307.5

# WVR row 91
$ws_WVR.Range("H91").Value = 24495.5
$ws_WVR.Range("J91").Value = 24495.5
$ws_WVR.Range("L91").Value = 24495.5
$ws_WVR.Range("N91").Value = -27303.5

# WVR row 100
$ws_WVR.Range("H100").Value = 6096.909
$ws_WVR.Range("I100").Value = 1344.3334
$ws_WVR.Range("K100").Value = 2688.6668
$ws_WVR.Range("M100").Value = -2147.6668

# WVR row 107
$ws_WVR.Range("H107").Value = 976.5
$ws_WVR.Range("I107").Value = 777.0714
$ws_WVR.Range("K107").Value = 2331.2142
$ws_WVR.Range("M107").Value = -411.2142000000003

# WVR row 122
$ws_WVR.Range("H122").Value = 2833.25
$ws_WVR.Range("I122").Value = 2305.5151
$ws_WVR.Range("K122").Value = 6916.5453
$ws_WVR.Range("M122").Value = -4466.5453

# WVR row 127
$ws_WVR.Range("H127").Value = 51832.832
$ws_WVR.Range("J127").Value = 58999.25
$ws_WVR.Range("L127").Value = 58999.25
$ws_WVR.Range("N127").Value = -68919.25

# WVR row 132
$ws_WVR.Range("H132").Value = 3843.25
$ws_WVR.Range("I132").Value = 3945.2
$ws_WVR.Range("K132").Value = 11835.6
$ws_WVR.Range("M132").Value = -9305.599999999999

# WVR row 136
$ws_WVR.Range("H136").Value = 5018.388
$ws_WVR.Range("I136").Value = 2897
$ws_WVR.Range("J136").Value = 9665.237999999999
$ws_WVR.Range("K136").Value = 8691
$ws_WVR.Range("L136").Value = 28995.714
$ws_WVR.Range("M136").Value = -6141
$ws_WVR.Range("N136").Value = -34095.714
